$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 12989582
$ws.Range("I32").Value2 = 22729020
$ws.Range("J32").Value2 = 3664.3333
$ws.Range("K32").Value2 = 22729020
$ws.Range("L32").Value2 = 3664.3333
$ws.Range("M32").Value2 = -22728694
$ws.Range("N32").Value2 = -4316.3333
$ws.Range("H111").Value2 = 3080.889
$ws.Range("I111").Value2 = 4553.5557
$ws.Range("J111").Value2 = 1608.2222
$ws.Range("K111").Value2 = 13660.6671
$ws.Range("L111").Value2 = 4824.6666
$ws.Range("M111").Value2 = -10593.6671
$ws.Range("N111").Value2 = -10958.6666
$ws.Range("H129").Value2 = 1734.0454
$ws.Range("J129").Value2 = 2593.1538
$ws.Range("L129").Value2 = 7779.4614
$ws.Range("N129").Value2 = -17779.4614
$ws.Range("H132").Value2 = 10582
$ws.Range("I132").Value2 = 11361.637
$ws.Range("K132").Value2 = 34084.911
$ws.Range("M132").Value2 = -31554.911
$ws.Range("H135").Value2 = 364.78946
$ws.Range("I135").Value2 = 330.3125
$ws.Range("J135").Value2 = 548.6667
$ws.Range("K135").Value2 = 2972.8125
$ws.Range("L135").Value2 = 4938.0003
$ws.Range("M135").Value2 = -437.8125
$ws.Range("N135").Value2 = -10008.0003

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2173544.8
$ws.Range("I32").Value2 = 1004005.44
$ws.Range("K32").Value2 = 1004005.44
$ws.Range("M32").Value2 = -1003718.44
$ws.Range("H74").Value2 = 1715.92
$ws.Range("I74").Value2 = 1339.3
$ws.Range("K74").Value2 = 1339.3
$ws.Range("M74").Value2 = -465.3
$ws.Range("H77").Value2 = 1715.92
$ws.Range("I77").Value2 = 1339.3
$ws.Range("K77").Value2 = 6696.5
$ws.Range("M77").Value2 = -2328.5
$ws.Range("H122").Value2 = 2336.7742
$ws.Range("I122").Value2 = 2405.3333
$ws.Range("J122").Value2 = 1874
$ws.Range("K122").Value2 = 7215.999899999999
$ws.Range("L122").Value2 = 5622
$ws.Range("M122").Value2 = -4765.999899999999
$ws.Range("N122").Value2 = -10522
$ws.Range("H125").Value2 = 87236
$ws.Range("J125").Value2 = 87236
$ws.Range("L125").Value2 = 87236
$ws.Range("N125").Value2 = -97076

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 3375.3333
$ws.Range("I99").Value2 = 2979.6667
$ws.Range("K99").Value2 = 2979.6667
$ws.Range("M99").Value2 = -1481.6667
$ws.Range("H134").Value2 = 3131.4546
$ws.Range("I134").Value2 = 2907.6667
$ws.Range("K134").Value2 = 8723.000100000001
$ws.Range("M134").Value2 = -6188.000100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1231.1111
$ws.Range("I16").Value2 = 1231.1111
$ws.Range("K16").Value2 = 1231.1111
$ws.Range("M16").Value2 = -944.1111000000001
$ws.Range("H31").Value2 = 2978962.2
$ws.Range("I31").Value2 = 2183.5625
$ws.Range("J31").Value2 = 12504654
$ws.Range("K31").Value2 = 2183.5625
$ws.Range("L31").Value2 = 12504654
$ws.Range("M31").Value2 = -1888.5625
$ws.Range("N31").Value2 = -12505244
$ws.Range("H34").Value2 = 2978962.2
$ws.Range("I34").Value2 = 2183.5625
$ws.Range("J34").Value2 = 12504654
$ws.Range("K34").Value2 = 2183.5625
$ws.Range("L34").Value2 = 12504654
$ws.Range("M34").Value2 = -1981.5625
$ws.Range("N34").Value2 = -12505058
$ws.Range("H58").Value2 = 2594.9
$ws.Range("I58").Value2 = 2049.75
$ws.Range("J58").Value2 = 2958.3333
$ws.Range("K58").Value2 = 2049.75
$ws.Range("L58").Value2 = 2958.3333
$ws.Range("M58").Value2 = -1846.75
$ws.Range("N58").Value2 = -3364.3333
$ws.Range("H99").Value2 = 6199.4
$ws.Range("I99").Value2 = 4999.5
$ws.Range("K99").Value2 = 4999.5
$ws.Range("M99").Value2 = -3501.5
$ws.Range("H113").Value2 = 1231.1111
$ws.Range("I113").Value2 = 1231.1111
$ws.Range("K113").Value2 = 1231.1111
$ws.Range("M113").Value2 = 938.8888999999999
$ws.Range("H122").Value2 = 4249.25
$ws.Range("J122").Value2 = 4499
$ws.Range("L122").Value2 = 13497
$ws.Range("N122").Value2 = -18397
$ws.Range("H126").Value2 = 6199.4
$ws.Range("I126").Value2 = 4999.5
$ws.Range("K126").Value2 = 14998.5
$ws.Range("M126").Value2 = -12528.5
$ws.Range("H132").Value2 = 4077.2258
$ws.Range("I132").Value2 = 3367.9092
$ws.Range("K132").Value2 = 10103.7276
$ws.Range("M132").Value2 = -7573.7276
$ws.Range("H134").Value2 = 4069.6072
$ws.Range("I134").Value2 = 4212.5
$ws.Range("J134").Value2 = 3212.25
$ws.Range("K134").Value2 = 12637.5
$ws.Range("L134").Value2 = 9636.75
$ws.Range("M134").Value2 = -10102.5
$ws.Range("N134").Value2 = -14706.75
$ws.Range("H136").Value2 = 2594.9
$ws.Range("I136").Value2 = 2049.75
$ws.Range("J136").Value2 = 2958.3333
$ws.Range("K136").Value2 = 6149.25
$ws.Range("L136").Value2 = 8874.999899999999
$ws.Range("M136").Value2 = -3599.25
$ws.Range("N136").Value2 = -13974.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value2 = 5365.385
$ws.Range("I54").Value2 = 0
$ws.Range("J54").Value2 = 5365.385
$ws.Range("K54").Value2 = 0
$ws.Range("L54").Value2 = 16096.155
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value2 = -17214.155

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 91820104
$ws.Range("I80").Value2 = 202000660
$ws.Range("J80").Value2 = 2974.6667
$ws.Range("K80").Value2 = 202000660
$ws.Range("L80").Value2 = 2974.6667
$ws.Range("M80").Value2 = -201999662
$ws.Range("N80").Value2 = -4970.6667
$ws.Range("H83").Value2 = 91820104
$ws.Range("I83").Value2 = 202000660
$ws.Range("J83").Value2 = 2974.6667
$ws.Range("K83").Value2 = 1010003300
$ws.Range("L83").Value2 = 14873.3335
$ws.Range("M83").Value2 = -1009998308
$ws.Range("N83").Value2 = -24857.3335
$ws.Range("H88").Value2 = 0
$ws.Range("I88").Value2 = 0
$ws.Range("K88").Value2 = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value2 = 0
$ws.Range("I91").Value2 = 0
$ws.Range("K91").Value2 = 0
$ws.Range("M91").ClearContents()
$ws.Range("H92").Value2 = 13966.667
$ws.Range("J92").Value2 = 13966.667
$ws.Range("L92").Value2 = 13966.667
$ws.Range("N92").Value2 = -17710.667
$ws.Range("H113").Value2 = 2145
$ws.Range("J113").Value2 = 2140.25
$ws.Range("L113").Value2 = 2140.25
$ws.Range("N113").Value2 = -6480.25
$ws.Range("H122").Value2 = 3560.04
$ws.Range("I122").Value2 = 2644.3333
$ws.Range("J122").Value2 = 4405.3076
$ws.Range("K122").Value2 = 7932.999899999999
$ws.Range("L122").Value2 = 13215.9228
$ws.Range("M122").Value2 = -5482.999899999999
$ws.Range("N122").Value2 = -18115.9228
$ws.Range("H132").Value2 = 2312.2222
$ws.Range("I132").Value2 = 1889
$ws.Range("K132").Value2 = 5667
$ws.Range("M132").Value2 = -3137
$ws.Range("H134").Value2 = 39667
$ws.Range("J134").Value2 = 39667
$ws.Range("L134").Value2 = 119001
$ws.Range("N134").Value2 = -124071

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2304.875
$ws.Range("I7").Value2 = 1906.5
$ws.Range("J7").Value2 = 3500
$ws.Range("K7").Value2 = 1906.5
$ws.Range("L7").Value2 = 3500
$ws.Range("M7").Value2 = -1794.5
$ws.Range("N7").Value2 = -3724
$ws.Range("H40").Value2 = 10399.9
$ws.Range("I40").Value2 = 10399.9
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 10399.9
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = -10263.9
$ws.Range("N40").ClearContents()
$ws.Range("H104").Value2 = 38000
$ws.Range("J104").Value2 = 38000
$ws.Range("L104").Value2 = 38000
$ws.Range("N104").Value2 = -44988
$ws.Range("H126").Value2 = 2304.875
$ws.Range("I126").Value2 = 1906.5
$ws.Range("J126").Value2 = 3500
$ws.Range("K126").Value2 = 5719.5
$ws.Range("L126").Value2 = 10500
$ws.Range("M126").Value2 = -3249.5
$ws.Range("N126").Value2 = -15440
$ws.Range("H132").Value2 = 3300.394
$ws.Range("I132").Value2 = 3302.8333
$ws.Range("K132").Value2 = 9908.499899999999
$ws.Range("M132").Value2 = -7378.499899999999
$ws.Range("H136").Value2 = 7812.875
$ws.Range("I136").Value2 = 7899.6
$ws.Range("K136").Value2 = 23698.8
$ws.Range("M136").Value2 = -21148.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value2 = 42208.082
$ws.Range("J54").Value2 = 42208.082
$ws.Range("L54").Value2 = 42208.082
$ws.Range("N54").Value2 = -43248.082
$ws.Range("H132").Value2 = 3154.0322
$ws.Range("I132").Value2 = 3213.3103
$ws.Range("J132").Value2 = 2294.5
$ws.Range("K132").Value2 = 9639.930899999999
$ws.Range("L132").Value2 = 6883.5
$ws.Range("M132").Value2 = -7109.930899999999
$ws.Range("N132").Value2 = -11943.5
